$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 10: "Meeting Management" across columns A-D
$ws.Range("A10").Value = "Meeting Management"
$ws.Range("B10").Value = "Meeting Management"
$ws.Range("C10").Value = "Meeting Management"
$ws.Range("D10").Value = "Meeting Management"

# Row 11: "View PALMS Summary" across columns A-D
$ws.Range("A11").Value = "View PALMS Summary"
$ws.Range("B11").Value = "View PALMS Summary"
$ws.Range("C11").Value = "View PALMS Summary"
$ws.Range("D11").Value = "View PALMS Summary"
